$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.000180381901313531
$ws.Range("B2").Value = 0.00202646505103288
$ws.Range("C2").Value = 0.00282783521043485
$ws.Range("D2").Value = 0.00298132524244602
$ws.Range("E2").Value = 0.0030861386483235
$ws.Range("F2").Value = 0.00305836028205187
$ws.Range("G2").Value = 0.00299568485409739
$ws.Range("H2").Value = 0.00288788323278704
$ws.Range("I2").Value = 0.0029699106583921
$ws.Range("J2").Value = 0.00393124807413163
$ws.Range("K2").Value = 0.00495386842668284
$ws.Range("L2").Value = 0.00630440130243923
$ws.Range("A3").Value = 0.000668304945679488
$ws.Range("B3").Value = 0.00179887328962021
$ws.Range("C3").Value = 0.00450723029236909
$ws.Range("D3").Value = 0.00577982155150165
$ws.Range("E3").Value = 0.00661684734541012
$ws.Range("F3").Value = 0.00688989949193351
$ws.Range("G3").Value = 0.00685754729176655
$ws.Range("H3").Value = 0.00656520033955415
$ws.Range("I3").Value = 0.0062000497547292
$ws.Range("J3").Value = 0.00601138061648165
$ws.Range("K3").Value = 0.00648197264228069
$ws.Range("L3").Value = 0.00710033568228292
$ws.Range("A4").Value = 0.00362097285396867
$ws.Range("B4").Value = 0.00669633822611369
$ws.Range("C4").Value = 0.00790639217994997
$ws.Range("D4").Value = 0.0086346801410224
$ws.Range("E4").Value = 0.00875732651170694
$ws.Range("F4").Value = 0.00858073486317134
$ws.Range("G4").Value = 0.00814887091239454
$ws.Range("H4").Value = 0.00762382637705408
$ws.Range("I4").Value = 0.00727776724187815
$ws.Range("J4").Value = 0.00755158189880095
$ws.Range("K4").Value = 0.00802179536531146
$ws.Range("L4").Value = 0.00843252399589791
$ws.Range("A5").Value = 0.00222208090089993
$ws.Range("B5").Value = 0.00198832872418692
$ws.Range("C5").Value = 0.00187962736755257
$ws.Range("D5").Value = 0.0016307441038607
$ws.Range("E5").Value = 0.00159585582949227
$ws.Range("F5").Value = 0.00229592863351751
$ws.Range("G5").Value = 0.00416891122526757
$ws.Range("H5").Value = 0.00588913972936315
$ws.Range("I5").Value = 0.00836865631176685
$ws.Range("J5").Value = 0.0103286146481506
$ws.Range("K5").Value = 0.0117492362698209
$ws.Range("L5").Value = 0.0123648123026585
$ws.Range("A6").Value = 0.00254775621877279
$ws.Range("B6").Value = 0.00326746152119006
$ws.Range("C6").Value = 0.00379947483457112
$ws.Range("D6").Value = 0.00404817982250153
$ws.Range("E6").Value = 0.00430491681395852
$ws.Range("F6").Value = 0.00493814770243131
$ws.Range("G6").Value = 0.00536901572958919
$ws.Range("H6").Value = 0.00650470712950334
$ws.Range("I6").Value = 0.00737992660640754
$ws.Range("J6").Value = 0.00792000404753205
$ws.Range("K6").Value = 0.00790966279725464
$ws.Range("L6").Value = 0.00759036333596649
$ws.Range("A7").Value = 0.00054319070712161
$ws.Range("B7").Value = 0.000497314088401135
$ws.Range("C7").Value = 0.00128007060736738
$ws.Range("D7").Value = 0.00282281644559549
$ws.Range("E7").Value = 0.00535307061564333
$ws.Range("F7").Value = 0.0074407486856341
$ws.Range("G7").Value = 0.0103275814820374
$ws.Range("H7").Value = 0.0125506899772611
$ws.Range("I7").Value = 0.0141464223719995
$ws.Range("J7").Value = 0.014858881927461
$ws.Range("K7").Value = 0.0148732249325642
$ws.Range("L7").Value = 0.0149188592574019
$ws.Range("A8").Value = 0.00156592066703887
$ws.Range("B8").Value = 0.00270845790906021
$ws.Range("C8").Value = 0.00444724413097699
$ws.Range("D8").Value = 0.00719209029255016
$ws.Range("E8").Value = 0.00942950593137015
$ws.Range("F8").Value = 0.0124477919766737
$ws.Range("G8").Value = 0.014763537492099
$ws.Range("H8").Value = 0.0164136051836574
$ws.Range("I8").Value = 0.0171110329389953
$ws.Range("J8").Value = 0.0170778516352685
$ws.Range("K8").Value = 0.0170902796804482
$ws.Range("L8").Value = 0.0170502601973408
$ws.Range("A9").Value = 0.00043885034746971
$ws.Range("B9").Value = 0.00183939460703719
$ws.Range("C9").Value = 0.00449527559480983
$ws.Range("D9").Value = 0.00637778259189354
$ws.Range("E9").Value = 0.0092396326888571
$ws.Range("F9").Value = 0.0112413897940532
$ws.Range("G9").Value = 0.012525109283382
$ws.Range("H9").Value = 0.0128307817181426
$ws.Range("I9").Value = 0.0124827074844336
$ws.Range("J9").Value = 0.012197329497849
$ws.Range("K9").Value = 0.0118900548865487
$ws.Range("L9").Value = 0.0117665838391631
$ws.Range("A10").Value = 0.00187190423219707
$ws.Range("B10").Value = 0.00487992664692627
$ws.Range("C10").Value = 0.00680670985297365
$ws.Range("D10").Value = 0.00994251575377175
$ws.Range("E10").Value = 0.0120813123490909
$ws.Range("F10").Value = 0.0134590533430708
$ws.Range("G10").Value = 0.0137972119855614
$ws.Range("H10").Value = 0.0134284447377566
$ws.Range("I10").Value = 0.0131660410075745
$ws.Range("J10").Value = 0.0128965887153549
$ws.Range("K10").Value = 0.0128710263140479
$ws.Range("L10").Value = 0.0132654339520382
$ws.Range("A11").Value = 0.00326073619670364
$ws.Range("B11").Value = 0.00442864254431445
$ws.Range("C11").Value = 0.00745540317367581
$ws.Range("D11").Value = 0.00927506466864859
$ws.Range("E11").Value = 0.0103279725968415
$ws.Range("F11").Value = 0.0103430428790442
$ws.Range("G11").Value = 0.00977526439376639
$ws.Range("H11").Value = 0.00934805094269594
$ws.Range("I11").Value = 0.00895633126316519
$ws.Range("J11").Value = 0.00879584949550368
$ws.Range("K11").Value = 0.00906284625012403
$ws.Range("L11").Value = 0.00956906819057357
$ws.Range("A12").Value = 0.000614631912208452
$ws.Range("B12").Value = 0.00280199229229439
$ws.Range("C12").Value = 0.00379585322807914
$ws.Range("D12").Value = 0.00424277821281329
$ws.Range("E12").Value = 0.00389566353102705
$ws.Range("F12").Value = 0.00382769585201196
$ws.Range("G12").Value = 0.00376663859902029
$ws.Range("H12").Value = 0.00384446597619
$ws.Range("I12").Value = 0.00369828284973274
$ws.Range("J12").Value = 0.00356103978318111
$ws.Range("K12").Value = 0.00363497225993592
$ws.Range("L12").Value = 0.00416878120389194
$ws.Range("A13").Value = 0.00485881825032308
$ws.Range("B13").Value = 0.00565197478074263
$ws.Range("C13").Value = 0.00593311115153079
$ws.Range("D13").Value = 0.00533485500451292
$ws.Range("E13").Value = 0.00491348713516027
$ws.Range("F13").Value = 0.00460720947783089
$ws.Range("G13").Value = 0.00447394023677111
$ws.Range("H13").Value = 0.00421602365009811
$ws.Range("I13").Value = 0.00405944761435302
$ws.Range("J13").Value = 0.00413670286088756
$ws.Range("K13").Value = 0.00468410349962178
$ws.Range("L13").Value = 0.00539009215060856
$ws.Range("A14").Value = 0.00276185259355888
$ws.Range("B14").Value = 0.00420391775995543
$ws.Range("C14").Value = 0.00689511046297617
$ws.Range("D14").Value = 0.0102476853899128
$ws.Range("E14").Value = 0.011907035085175
$ws.Range("F14").Value = 0.0132275892111421
$ws.Range("G14").Value = 0.0137872377364079
$ws.Range("H14").Value = 0.013725378151362
$ws.Range("I14").Value = 0.0134884678803199
$ws.Range("J14").Value = 0.0130696641458353
$ws.Range("K14").Value = 0.0126362546007619
$ws.Range("L14").Value = 0.0121601820369854
$ws.Range("A15").Value = 0.000181477450783341
$ws.Range("B15").Value = 0.00253871697891218
$ws.Range("C15").Value = 0.00562880214420549
$ws.Range("D15").Value = 0.00666112731115839
$ws.Range("E15").Value = 0.0075091245274166
$ws.Range("F15").Value = 0.0076379682358312
$ws.Range("G15").Value = 0.00726610956218062
$ws.Range("H15").Value = 0.0068442537168651
$ws.Range("I15").Value = 0.00646418318748947
$ws.Range("J15").Value = 0.00620479019597399
$ws.Range("K15").Value = 0.00618129041533033
$ws.Range("A16").Value = 0.00356062497141529
$ws.Range("B16").Value = 0.00690756920212721
$ws.Range("C16").Value = 0.0077911950172133
$ws.Range("D16").Value = 0.0085862629009385
$ws.Range("E16").Value = 0.00864528144771831
$ws.Range("F16").Value = 0.00818390322065299
$ws.Range("G16").Value = 0.00768571502177554
$ws.Range("H16").Value = 0.00718966277794134
$ws.Range("I16").Value = 0.00679407783334296
$ws.Range("J16").Value = 0.00657778430539083
$ws.Range("A17").Value = 0.00243179312833242
$ws.Range("B17").Value = 0.00178087060698805
$ws.Range("C17").Value = 0.0015434305732571
$ws.Range("D17").Value = 0.00153128225928245
$ws.Range("E17").Value = 0.00268211989224787
$ws.Range("F17").Value = 0.003854854698073
$ws.Range("G17").Value = 0.00541364606121893
$ws.Range("H17").Value = 0.00673658854476051
$ws.Range("I17").Value = 0.00820883364744161
$ws.Range("A18").Value = 0.00363653410277287
$ws.Range("B18").Value = 0.0040275990369249
$ws.Range("C18").Value = 0.00524251775916034
$ws.Range("D18").Value = 0.00712203330186178
$ws.Range("E18").Value = 0.00864918276896963
$ws.Range("F18").Value = 0.0103806357261635
$ws.Range("G18").Value = 0.0118372814920623
$ws.Range("H18").Value = 0.0133890355727637
$ws.Range("A19").Value = 0.00229038801941428
$ws.Range("B19").Value = 0.00186457041595478
$ws.Range("C19").Value = 0.00187600870927376
$ws.Range("D19").Value = 0.00249602363031546
$ws.Range("E19").Value = 0.00390217264997264
$ws.Range("F19").Value = 0.00510556247736726
$ws.Range("G19").Value = 0.0065338349648175
$ws.Range("A20").Value = 0.00289287928538506
$ws.Range("B20").Value = 0.00547995357060707
$ws.Range("C20").Value = 0.00710931106184179
$ws.Range("D20").Value = 0.00901218328811792
$ws.Range("E20").Value = 0.0105009786811724
$ws.Range("F20").Value = 0.0121014854402187
$ws.Range("A21").Value = 0.00193271564719311
$ws.Range("B21").Value = 0.00249555092316173
$ws.Range("C21").Value = 0.00395846175475973
$ws.Range("D21").Value = 0.00506515893075026
$ws.Range("E21").Value = 0.00646325464058241
$ws.Range("A22").Value = 0.000460042284140094
$ws.Range("B22").Value = 0.0012752926844503
$ws.Range("C22").Value = 0.00204335883363688
$ws.Range("D22").Value = 0.00344405149113161
$ws.Range("A23").Value = 0.00252358987533086
$ws.Range("B23").Value = 0.00338151817296153
$ws.Range("C23").Value = 0.00487690211222989
$ws.Range("A24").Value = 0.000619765952155404
$ws.Range("B24").Value = 0.000905708615045294
$ws.Range("A25").Value = 0.00220995805780966
